# edits in upload timesheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record 2 sick days ("مرضى") for the first employee row (F2)
$ws.Range("F2").Value = 2

# Move the selection/active cell to F10 (matches the saved view state)
$ws.Range("F10").Select()
